# Scheduled market-price refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) per sheet with latest Universalis data. Generated from the nightly diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 233.55556
$ws.Range("I2").Value = 241.66667
$ws.Range("J2").Value = 217.33333
$ws.Range("K2").Value = 241.66667
$ws.Range("L2").Value = 217.33333
$ws.Range("M2").Value = -128.66667
$ws.Range("N2").Value = -443.33333

# Row 69
$ws.Range("H69").Value = 3314.2856
$ws.Range("I69").Value = 1733.3334
$ws.Range("K69").Value = 5200.0002
$ws.Range("M69").Value = -4326.0002

# Row 72
$ws.Range("H72").Value = 3314.2856
$ws.Range("I72").Value = 1733.3334
$ws.Range("K72").Value = 15600.0006
$ws.Range("M72").Value = -11232.0006

# Row 138
$ws.Range("H138").Value = 1726627.5
$ws.Range("I138").Value = 2704378.8
$ws.Range("J138").Value = 3922.8096
$ws.Range("K138").Value = 8113136.399999999
$ws.Range("L138").Value = 11768.4288
$ws.Range("M138").Value = -8107996.399999999
$ws.Range("N138").Value = -22048.4288

# Row 141
$ws.Range("H141").Value = 1606.826
$ws.Range("I141").Value = 807.6429000000001
$ws.Range("J141").Value = 2850
$ws.Range("K141").Value = 2422.9287
$ws.Range("L141").Value = 8550
$ws.Range("M141").Value = 2757.0713
$ws.Range("N141").Value = -18910

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5160.017
$ws.Range("I32").Value = 5427.479
$ws.Range("K32").Value = 5427.479
$ws.Range("M32").Value = -5140.479

# Row 61
$ws.Range("H61").Value = 1589.04
$ws.Range("I61").Value = 1606.2222
$ws.Range("J61").Value = 1544.8572
$ws.Range("K61").Value = 1606.2222
$ws.Range("L61").Value = 1544.8572
$ws.Range("M61").Value = -1394.2222
$ws.Range("N61").Value = -1968.8572

# Row 97
$ws.Range("H97").Value = 491.42856
$ws.Range("I97").Value = 370.9091
$ws.Range("K97").Value = 370.9091
$ws.Range("M97").Value = 125.0909

# Row 102
$ws.Range("H102").Value = 1746.3636
$ws.Range("I102").Value = 1746.3636
$ws.Range("K102").Value = 1746.3636
$ws.Range("M102").Value = -124.3635999999999

# Row 122
$ws.Range("H122").Value = 1001.8095
$ws.Range("I122").Value = 948.58826
$ws.Range("J122").Value = 1228
$ws.Range("K122").Value = 2845.76478
$ws.Range("L122").Value = 3684
$ws.Range("M122").Value = -395.76478
$ws.Range("N122").Value = -8584

# Row 132
$ws.Range("H132").Value = 1733.1372
$ws.Range("I132").Value = 1642.5682
$ws.Range("J132").Value = 2302.4285
$ws.Range("K132").Value = 4927.7046
$ws.Range("L132").Value = 6907.2855
$ws.Range("M132").Value = -2397.7046
$ws.Range("N132").Value = -11967.2855

# Row 136
$ws.Range("H136").Value = 1589.04
$ws.Range("I136").Value = 1606.2222
$ws.Range("J136").Value = 1544.8572
$ws.Range("K136").Value = 4818.6666
$ws.Range("L136").Value = 4634.571599999999
$ws.Range("M136").Value = -2268.6666
$ws.Range("N136").Value = -9734.571599999999

$ws = $wb.Worksheets.Item("BSM")
# Row 27
$ws.Range("H27").Value = 35000
$ws.Range("J27").Value = 35000
$ws.Range("L27").Value = 35000
$ws.Range("N27").Value = -35384

# Row 134
$ws.Range("H134").Value = 30622.234
$ws.Range("I134").Value = 35623.69
$ws.Range("J134").Value = 1613.8
$ws.Range("K134").Value = 106871.07
$ws.Range("L134").Value = 4841.4
$ws.Range("M134").Value = -104336.07
$ws.Range("N134").Value = -9911.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2103.25
$ws.Range("I31").Value = 1894
$ws.Range("K31").Value = 1894
$ws.Range("M31").Value = -1599

# Row 34
$ws.Range("H34").Value = 2103.25
$ws.Range("I34").Value = 1894
$ws.Range("K34").Value = 1894
$ws.Range("M34").Value = -1692

# Row 74
$ws.Range("H74").Value = 14460.5
$ws.Range("J74").Value = 14460.5
$ws.Range("L74").Value = 14460.5
$ws.Range("N74").Value = -16208.5

# Row 77
$ws.Range("H77").Value = 14460.5
$ws.Range("J77").Value = 14460.5
$ws.Range("L77").Value = 43381.5
$ws.Range("N77").Value = -52117.5

# Row 103
$ws.Range("H103").Value = 29351.715
$ws.Range("I103").Value = 10000
$ws.Range("J103").Value = 32577
$ws.Range("K103").Value = 10000
$ws.Range("L103").Value = 32577
$ws.Range("M103").Value = -8828
$ws.Range("N103").Value = -34921

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null

# Row 132
$ws.Range("H132").Value = 1722.381
$ws.Range("I132").Value = 1615.973
$ws.Range("J132").Value = 2509.8
$ws.Range("K132").Value = 4847.919
$ws.Range("L132").Value = 7529.400000000001
$ws.Range("M132").Value = -2317.919
$ws.Range("N132").Value = -12589.4

# Row 134
$ws.Range("H134").Value = 4659.1665
$ws.Range("I134").Value = 5314
$ws.Range("J134").Value = 1385
$ws.Range("K134").Value = 15942
$ws.Range("L134").Value = 4155
$ws.Range("M134").Value = -13407
$ws.Range("N134").Value = -9225

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 45515.953
$ws.Range("I12").Value = 75.125
$ws.Range("J12").Value = 71482.14
$ws.Range("K12").Value = 225.375
$ws.Range("L12").Value = 214446.42
$ws.Range("M12").Value = -52.375
$ws.Range("N12").Value = -214792.42

# Row 64
$ws.Range("H64").Value = 1078341
$ws.Range("I64").Value = 308.25
$ws.Range("J64").Value = 1557466.6
$ws.Range("K64").Value = 924.75
$ws.Range("L64").Value = 4672399.800000001
$ws.Range("M64").Value = -654.75
$ws.Range("N64").Value = -4672939.800000001

# Row 67
$ws.Range("H67").Value = 1078341
$ws.Range("I67").Value = 308.25
$ws.Range("J67").Value = 1557466.6
$ws.Range("K67").Value = 924.75
$ws.Range("L67").Value = 4672399.800000001
$ws.Range("M67").Value = 11.25
$ws.Range("N67").Value = -4674271.800000001

# Row 133
$ws.Range("H133").Value = 9200
$ws.Range("I133").Value = 6000
$ws.Range("K133").Value = 18000
$ws.Range("M133").Value = -12940

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2512.9312
$ws.Range("I122").Value = 1326.7
$ws.Range("J122").Value = 5149
$ws.Range("K122").Value = 3980.1
$ws.Range("L122").Value = 15447
$ws.Range("M122").Value = -1530.1
$ws.Range("N122").Value = -20347

# Row 132
$ws.Range("H132").Value = 2716.9412
$ws.Range("I132").Value = 2353.077
$ws.Range("J132").Value = 3899.5
$ws.Range("K132").Value = 7059.231000000001
$ws.Range("L132").Value = 11698.5
$ws.Range("M132").Value = -4529.231000000001
$ws.Range("N132").Value = -16758.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1200.1052
$ws.Range("I46").Value = 1133.3334
$ws.Range("J46").Value = 1314.5714
$ws.Range("K46").Value = 1133.3334
$ws.Range("L46").Value = 1314.5714
$ws.Range("M46").Value = -945.3334
$ws.Range("N46").Value = -1690.5714

# Row 68
$ws.Range("H68").Value = 1767.2858
$ws.Range("I68").Value = 1683.7778
$ws.Range("J68").Value = 2049.125
$ws.Range("K68").Value = 1683.7778
$ws.Range("L68").Value = 2049.125
$ws.Range("M68").Value = -934.7778000000001
$ws.Range("N68").Value = -3547.125

# Row 71
$ws.Range("H71").Value = 1767.2858
$ws.Range("I71").Value = 1683.7778
$ws.Range("J71").Value = 2049.125
$ws.Range("K71").Value = 8418.889000000001
$ws.Range("L71").Value = 10245.625
$ws.Range("M71").Value = -4674.889000000001
$ws.Range("N71").Value = -17733.625

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4841.1763
$ws.Range("J81").Value = 5226.6665
$ws.Range("L81").Value = 10453.333
$ws.Range("N81").Value = -12575.333

# Row 84
$ws.Range("H84").Value = 4841.1763
$ws.Range("J84").Value = 5226.6665
$ws.Range("L84").Value = 52266.665
$ws.Range("N84").Value = -62874.665
